$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format to preserve exact string representation
# (avoids Excel auto-converting numeric-looking strings to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.763.13'
$ws.Range("E2").Value = '  +1.61%  '

$ws.Range("D3").Value = '2.105.47'
$ws.Range("E3").Value = '  +5.30%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = '333.75'
$ws.Range("E5").Value = '  +2.96%  '

$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("D7").Value = '0.5295'
$ws.Range("E7").Value = '  +3.96%  '

$ws.Range("D8").Value = '0.4370'
$ws.Range("E8").Value = '  +5.58%  '

$ws.Range("D9").Value = '0.08939'
$ws.Range("E9").Value = '  +2.72%  '

$ws.Range("D10").Value = '46.92'
$ws.Range("E10").Value = '  +9.47%  '

$ws.Range("D11").Value = '1.166'
$ws.Range("E11").Value = '  +2.79%  '

$ws.Range("D12").Value = '24.77'
$ws.Range("E12").Value = '  -0.86%  '

$ws.Range("D13").Value = '2.105.78'
$ws.Range("E13").Value = '  +5.00%  '

$ws.Range("D14").Value = '6.725'
$ws.Range("E14").Value = '  +3.09%  '

$ws.Range("D15").Value = '7.760'
$ws.Range("E15").Value = '  +4.48%  '

$ws.Range("D16").Value = '96.61'
$ws.Range("E16").Value = '  +2.58%  '

$ws.Range("D17").Value = '1.004'
$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("D18").Value = '0.00001131'
$ws.Range("E18").Value = '  +1.35%  '

$ws.Range("D19").Value = '0.06686'
$ws.Range("E19").Value = '  +2.13%  '

$ws.Range("D20").Value = '19.02'
$ws.Range("E20").Value = '  +0.52%  '

$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("D22").Value = '6.301'
$ws.Range("E22").Value = '  +2.57%  '

$ws.Range("D23").Value = '30.828.41'

$ws.Range("D24").Value = '12.24'
$ws.Range("E24").Value = '  +4.62%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.352.34'
$ws.Range("E25").Value = '  +5.08%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '2.289'
$ws.Range("E26").Value = '  +3.50%  '

$ws.Range("D27").Value = '22.60'
$ws.Range("E27").Value = '  +0.19%  '

$ws.Range("D28").Value = '2.574'
$ws.Range("E28").Value = '  +7.63%  '

$ws.Range("D29").Value = '162.73'
$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("D30").Value = '132.86'
$ws.Range("E30").Value = '  +0.92%  '

$ws.Range("D31").Value = '1.190'
$ws.Range("E31").Value = '  +4.50%  '

$ws.Range("D32").Value = '0.1080'
$ws.Range("E32").Value = '  +2.55%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.169'
$ws.Range("E33").Value = '  +1.79%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '4.035'
$ws.Range("E34").Value = '  +5.70%  '

$ws.Range("D35").Value = '1.539'
$ws.Range("E35").Value = '  +14.70%  '

$ws.Range("D36").Value = '0.02598'
$ws.Range("E36").Value = '  +3.81%  '

$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '9.574'
$ws.Range("E37").Value = '  +7.46%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '5.524'
$ws.Range("E38").Value = '  +2.40%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.06742'
$ws.Range("E39").Value = '  +2.61%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '12.64'
$ws.Range("E40").Value = '  +3.16%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2271'
$ws.Range("E41").Value = '  +3.18%  '

$ws.Range("D42").Value = '0.6802'
$ws.Range("E42").Value = '  +3.11%  '

$ws.Range("D43").Value = '1.245'
$ws.Range("E43").Value = '  +1.84%  '

$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.15%  '

$ws.Range("D45").Value = '14.07'
$ws.Range("E45").Value = '  +3.29%  '

$ws.Range("D46").Value = '0.6403'
$ws.Range("E46").Value = '  +4.33%  '

$ws.Range("D47").Value = '2.217'
$ws.Range("E47").Value = '  +0.17%  '

$ws.Range("D48").Value = '3.659'
$ws.Range("E48").Value = '  -0.18%  '

$ws.Range("D49").Value = '1.257'
$ws.Range("E49").Value = '  -0.10%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '82.90'
$ws.Range("E50").Value = '  +3.67%  '

$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Value = '1.190'
$ws.Range("E51").Value = '  +7.89%  '

# Restore default style on column D (remove temporary text formatting)
$ws.Range("D2:D51").Style = "Normal"